# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" and before the
#    existing "2022-Q3" sheet, by duplicating the "2022-Q3" sheet (same
#    columns/headers/fund row) and updating the quarter-specific figures.
# 2) Update the "总计" (total) summary sheet: add a new first data row for
#    2022-Q4 and push the existing quarters (2022-Q3, 2022-Q2, 2021-Q4,
#    2020-Q4) down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the "2022-Q3" sheet to become the new "2022-Q4" sheet.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $q3)

$copied = $wb.Worksheets.Item(3)
$copied.Name = "2022-Q4"
$copied.Move($wb.Worksheets.Item("2022-Q3"))

# Re-fetch by name (position-bound references can go stale after Move).
$q4 = $wb.Worksheets.Item("2022-Q4")

# Update the quarter-specific numbers on the new sheet (fund stays the
# same: 005167 / 嘉实润泽量化一年定期开放混合). D:G keep their original
# text formatting; H2 is a real number.
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "0.56"
$q4.Range("E2").Value = "27.25"
$q4.Range("F2").Value = "0.55"
$q4.Range("G2").Value = "0.0031"
$q4.Range("H2").Value = 10

# ---------------------------------------------------------------------
# Step 2: shift the "总计" rows down to make room for 2022-Q4 at row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Give the new row 6 (col A) the same style as the other index cells
# (A2:A5) before writing into it.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

$total.Range("A6").Value = 4
$total.Range("B6").Value = $total.Range("B5").Value2
$total.Range("C6").Value = $total.Range("C5").Value2
$total.Range("D6").Value = $total.Range("D5").Value2

$total.Range("A5").Value = 3
$total.Range("B5").Value = $total.Range("B4").Value2
$total.Range("C5").Value = $total.Range("C4").Value2
$total.Range("D5").Value = $total.Range("D4").Value2

$total.Range("A4").Value = 2
$total.Range("B4").Value = $total.Range("B3").Value2
$total.Range("C4").Value = $total.Range("C3").Value2
$total.Range("D4").Value = $total.Range("D3").Value2

$total.Range("A3").Value = 1
$total.Range("B3").Value = $total.Range("B2").Value2
$total.Range("C3").Value = $total.Range("C2").Value2
$total.Range("D3").Value = $total.Range("D2").Value2

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0

# Restore the originally-active tab (the last sheet, "2020-Q4") since
# sheet copy/move operations above shift which tab is marked active.
$wb.Worksheets.Item("2020-Q4").Activate()
